# Trade #25 closed at 2026-02-17 12:37:35 - unknown UNKNOWN +0.000%
$wb = $excel.ActiveWorkbook

# --- Summary sheet updates ---
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B5").Value = 0.55    # Total P&L %
$summary.Range("B6").Value = 25      # Total Trades
$summary.Range("B9").Value = 40      # Win Rate %

# --- Strategy Status sheet updates (MarketMaking row) ---
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 25       # Trades
$status.Range("G4").Value = 40       # Win Rate %

# --- Append the new trade (#25) row to both the "All Trades" log and the
#     per-strategy "MarketMaking" log. Both sheets share identical layout. ---
function Add-TradeRow($ws) {
    $row = 26

    # Date / Time are text-like strings that Excel would otherwise auto-
    # convert to date/time serials, so force the cell to Text format first,
    # then restore the default "Normal" style so no stray style is left
    # behind on the cell.
    $ws.Cells.Item($row, 1).Value = 25

    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 2).Style = "Normal"

    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "12:37:29"
    $ws.Cells.Item($row, 3).Style = "Normal"

    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.32
    $ws.Cells.Item($row, 7).Value = 0.32
    $ws.Cells.Item($row, 8).Value = "CLOSED"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 100.69
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"
    $ws.Cells.Item($row, 16).Value = "early_exit"
    $ws.Cells.Item($row, 17).Value = 0.13
}

Add-TradeRow $wb.Worksheets.Item("All Trades")
Add-TradeRow $wb.Worksheets.Item("MarketMaking")
